# Update column F ("dSF") values on Sheet1 to match the re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -4
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = -5
$ws.Range("F13").Value = -4
$ws.Range("F14").Value = -8
$ws.Range("F16").Value = -7
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 4
$ws.Range("F20").Value = -2
